$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep Text format so numeric-looking strings
# (e.g. "315.21", "1.680", "27.376.65") are not coerced into numbers,
# which would lose formatting like trailing zeros or thousand-separator dots.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.376.65'
$ws.Range("E2").Value = '  +1.96%  '
$ws.Range("D3").Value = '1.840.80'
$ws.Range("E3").Value = '  +1.54%  '
$ws.Range("E4").Value = '  +1.31%  '
$ws.Range("D5").Value = '315.21'
$ws.Range("E5").Value = '  +2.06%  '
$ws.Range("D6").Value = '1.012'
$ws.Range("E6").Value = '  +1.14%  '
$ws.Range("E7").Value = '  +1.92%  '
$ws.Range("D8").Value = '0.3700'
$ws.Range("E8").Value = '  +0.48%  '
$ws.Range("D9").Value = '0.07471'
$ws.Range("E9").Value = '  +1.49%  '
$ws.Range("D10").Value = '0.8857'
$ws.Range("E10").Value = '  +1.84%  '
$ws.Range("D11").Value = '20.52'
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("D12").Value = '1.844.54'
$ws.Range("E12").Value = '  +1.29%  '
$ws.Range("D13").Value = '0.07393'
$ws.Range("E13").Value = '  +4.66%  '
$ws.Range("D14").Value = '5.485'
$ws.Range("E14").Value = '  +2.37%  '
$ws.Range("D15").Value = '93.33'
$ws.Range("E15").Value = '  +2.02%  '
$ws.Range("D16").Value = '6.579'
$ws.Range("E16").Value = '  +1.26%  '
$ws.Range("E17").Value = '  +1.16%  '
$ws.Range("D18").Value = '0.000008854'
$ws.Range("E18").Value = '  +1.94%  '
$ws.Range("D19").Value = '1.013'
$ws.Range("E19").Value = '  +1.23%  '
$ws.Range("D20").Value = '14.85'
$ws.Range("E20").Value = '  +0.87%  '
$ws.Range("D21").Value = '27.383.83'
$ws.Range("E21").Value = '  +1.86%  '
$ws.Range("D22").Value = '5.360'
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("E23").Value = '  +1.49%  '
$ws.Range("D24").Value = '2.075.70'
$ws.Range("E24").Value = '  +1.03%  '
$ws.Range("E25").Value = '  +0.53%  '
$ws.Range("D26").Value = '152.07'
$ws.Range("E26").Value = '  +1.09%  '
$ws.Range("D27").Value = '18.67'
$ws.Range("E27").Value = '  +2.03%  '
$ws.Range("D28").Value = '2.186'
$ws.Range("E28").Value = '  +0.48%  '
$ws.Range("D29").Value = '5.280'
$ws.Range("E29").Value = '  -0.60%  '
$ws.Range("D30").Value = '118.06'
$ws.Range("E30").Value = '  +2.17%  '
$ws.Range("E31").Value = '  +0.43%  '
$ws.Range("D32").Value = '0.7623'
$ws.Range("E32").Value = '  -0.41%  '
$ws.Range("D33").Value = '1.181'
$ws.Range("E34").Value = '  +1.32%  '
$ws.Range("D35").Value = '2.940'
$ws.Range("E35").Value = '  +1.33%  '
$ws.Range("E36").Value = '  +1.26%  '
$ws.Range("E37").Value = '  +1.81%  '
$ws.Range("E38").Value = '  +1.76%  '
$ws.Range("D39").Value = '0.01964'
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("D40").Value = '3.002'
$ws.Range("E40").Value = '  +2.11%  '
$ws.Range("D41").Value = '7.316'
$ws.Range("E41").Value = '  +0.90%  '
$ws.Range("D42").Value = '0.5360'
$ws.Range("E42").Value = '  +0.89%  '
$ws.Range("E43").Value = '  +1.57%  '
$ws.Range("D44").Value = '0.1669'
$ws.Range("E44").Value = '  +0.57%  '
$ws.Range("D45").Value = '8.555'
$ws.Range("E45").Value = '  +1.67%  '
$ws.Range("D46").Value = '0.4990'
$ws.Range("E46").Value = '  +1.39%  '
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("D48").Value = '1.014'
$ws.Range("D49").Value = '105.17'
$ws.Range("E49").Value = '  +1.50%  '
$ws.Range("D50").Value = '1.680'
$ws.Range("E50").Value = '  +0.64%  '
$ws.Range("D51").Value = '0.06334'
$ws.Range("E51").Value = '  +0.77%  '
